$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.775553848681341
$ws.Range("D2").Value = 3.758242158487112
$ws.Range("E2").Value = 16.5821330812193
$ws.Range("F2").Value = 19.2790344371383
$ws.Range("G2").Value = 20.41577253843935
$ws.Range("H2").Value = 11.97370744959122
$ws.Range("K2").Value = 11.53952637577816
$ws.Range("N2").Value = 17.13445589548761
$ws.Range("O2").Value = 17.12819954464804

$ws.Range("B3").Value = 7.701920712875018
$ws.Range("D3").Value = 3.685205191049644
$ws.Range("E3").Value = 15.63701532668809
$ws.Range("F3").Value = 19.23344622197329
$ws.Range("G3").Value = 20.29293657291592
$ws.Range("H3").Value = 12.00830172688339
$ws.Range("K3").Value = 11.10627725588193
$ws.Range("N3").Value = 17.18192543703669
$ws.Range("O3").Value = 17.15705391767574

$ws.Range("B4").Value = 7.658180701587655
$ws.Range("D4").Value = 3.639026346392571
$ws.Range("E4").Value = 15.03125312752327
$ws.Range("F4").Value = 19.21176221086434
$ws.Range("G4").Value = 20.22621151936711
$ws.Range("H4").Value = 12.03183791918186
$ws.Range("K4").Value = 10.82948156748783
$ws.Range("N4").Value = 17.21283981252195
$ws.Range("O4").Value = 17.17951307008911

$ws.Range("B5").Value = 7.640746443758189
$ws.Range("D5").Value = 3.619886478208135
$ws.Range("E5").Value = 14.77826600100212
$ws.Range("F5").Value = 19.20451750047833
$ws.Range("G5").Value = 20.20123208975739
$ws.Range("H5").Value = 12.04200526744054
$ws.Range("K5").Value = 10.71408110500883
$ws.Range("N5").Value = 17.22588314129175
$ws.Range("O5").Value = 17.18985395883944

$ws.Range("B6").Value = 7.637875640258342
$ws.Range("D6").Value = 3.616689309914457
$ws.Range("E6").Value = 14.73589604746348
$ws.Range("F6").Value = 19.20341078323142
$ws.Range("G6").Value = 20.19721847564928
$ws.Range("H6").Value = 12.04372831708821
$ws.Range("K6").Value = 10.69476504367038
$ws.Range("N6").Value = 17.22807590867906
$ws.Range("O6").Value = 17.19164272006527

$ws.Range("B7").Value = 7.657943971156499
$ws.Range("D7").Value = 3.638769502911935
$ws.Range("E7").Value = 15.02786569967539
$ws.Range("F7").Value = 19.21165805546455
$ws.Range("G7").Value = 20.22586565446149
$ws.Range("H7").Value = 12.0319727080954
$ws.Range("K7").Value = 10.82793562950414
$ws.Range("N7").Value = 17.2130139145013
$ws.Range("O7").Value = 17.17964772391144

$ws.Range("B8").Value = 7.749872192308528
$ws.Range("D8").Value = 3.733343904613907
$ws.Range("E8").Value = 16.26168442563297
$ws.Range("F8").Value = 19.2620105330386
$ws.Range("G8").Value = 20.37162998210168
$ws.Range("H8").Value = 11.98515861675695
$ws.Range("K8").Value = 11.39243913903358
$ws.Range("N8").Value = 17.15045702090274
$ws.Range("O8").Value = 17.13716201669867

$ws.Range("B9").Value = 7.940893604539491
$ws.Range("D9").Value = 3.907652333583437
$ws.Range("E9").Value = 18.58296039659888
$ws.Range("F9").Value = 19.41045643064052
$ws.Range("G9").Value = 20.72516011178071
$ws.Range("H9").Value = 11.91161015465655
$ws.Range("K9").Value = 12.40986523376614
$ws.Range("N9").Value = 17.04176911541525
$ws.Range("O9").Value = 17.09163205230656

$ws.Range("B10").Value = 8.086521371934836
$ws.Range("D10").Value = 4.028211690046437
$ws.Range("E10").Value = 20.23237390283409
$ws.Range("F10").Value = 19.54923241407178
$ws.Range("G10").Value = 21.02408385219091
$ws.Range("H10").Value = 11.86875889186055
$ws.Range("K10").Value = 13.09817664374714
$ws.Range("N10").Value = 16.9703852590758
$ws.Range("O10").Value = 17.08139703984631

$ws.Range("B11").Value = 8.153642599186256
$ws.Range("D11").Value = 4.081293269232784
$ws.Range("E11").Value = 20.94019214055777
$ws.Range("F11").Value = 19.61865100670179
$ws.Range("G11").Value = 21.16805413211683
$ws.Range("H11").Value = 11.85170475613829
$ws.Range("K11").Value = 13.39769036683442
$ws.Range("N11").Value = 16.93973803132141
$ws.Range("O11").Value = 17.08180954876361

$ws.Range("B12").Value = 8.179162161196688
$ws.Range("D12").Value = 4.101130330921731
$ws.Range("E12").Value = 21.20214799134854
$ws.Range("F12").Value = 19.64582544905518
$ws.Range("G12").Value = 21.22367146533613
$ws.Range("H12").Value = 11.84559840030725
$ws.Range("K12").Value = 13.50910022940463
$ws.Range("N12").Value = 16.9283943827463
$ws.Range("O12").Value = 17.08269585908289

$ws.Range("B13").Value = 8.17366189443665
$ws.Range("D13").Value = 4.096869955324108
$ws.Range("E13").Value = 21.14600097742708
$ws.Range("F13").Value = 19.63993378359989
$ws.Range("G13").Value = 21.21164522622802
$ws.Range("H13").Value = 11.84689785816247
$ws.Range("K13").Value = 13.48519626324113
$ws.Range("N13").Value = 16.93082580867622
$ws.Range("O13").Value = 17.08247249117226

$ws.Range("B14").Value = 8.155740194823551
$ws.Range("D14").Value = 4.082930620142458
$ws.Range("E14").Value = 20.96186506074546
$ws.Range("F14").Value = 19.6208689897081
$ws.Range("G14").Value = 21.1726080751599
$ws.Range("H14").Value = 11.8511953293247
$ws.Range("K14").Value = 13.40689667487902
$ws.Range("N14").Value = 16.93879953997475
$ws.Range("O14").Value = 17.08186782864489

$ws.Range("B15").Value = 8.144775250689449
$ws.Range("D15").Value = 4.074357724680347
$ws.Range("E15").Value = 20.84828573525881
$ws.Range("F15").Value = 19.60930624329327
$ws.Range("G15").Value = 21.14883827666905
$ws.Range("H15").Value = 11.85387347975193
$ws.Range("K15").Value = 13.35867276885017
$ws.Range("N15").Value = 16.94371775213388
$ws.Range("O15").Value = 17.0815925633303

$ws.Range("B16").Value = 8.082150423812138
$ws.Range("D16").Value = 4.024706325248706
$ws.Range("E16").Value = 20.18526377498149
$ws.Range("F16").Value = 19.54482086248701
$ws.Range("G16").Value = 21.01483161181327
$ws.Range("H16").Value = 11.86992258202627
$ws.Range("K16").Value = 13.07832424924664
$ws.Range("N16").Value = 16.97242480095802
$ws.Range("O16").Value = 17.0814721937267

$ws.Range("B17").Value = 8.043939313829023
$ws.Range("D17").Value = 3.9937879444143
$ws.Range("E17").Value = 19.76765750915807
$ws.Range("F17").Value = 19.50686020116409
$ws.Range("G17").Value = 20.93463438972751
$ws.Range("H17").Value = 11.88039358415539
$ws.Range("K17").Value = 12.90281499546255
$ws.Range("N17").Value = 16.99050270900521
$ws.Range("O17").Value = 17.08269759755239

$ws.Range("B18").Value = 8.02204500554039
$ws.Range("D18").Value = 3.975839289081588
$ws.Range("E18").Value = 19.52345839110827
$ws.Range("F18").Value = 19.48561929929227
$ws.Range("G18").Value = 20.88926215551418
$ws.Range("H18").Value = 11.88664575649735
$ws.Range("K18").Value = 12.80058903493372
$ws.Range("N18").Value = 17.0010725247304
$ws.Range("O18").Value = 17.08387944801963

$ws.Range("B19").Value = 8.014647064932912
$ws.Range("D19").Value = 3.969734137846628
$ws.Range("E19").Value = 19.44008822797822
$ws.Range("F19").Value = 19.4785298489601
$ws.Range("G19").Value = 20.87403105252487
$ws.Range("H19").Value = 11.88880202714657
$ws.Range("K19").Value = 12.76575939727473
$ws.Range("N19").Value = 17.00468082617223
$ws.Range("O19").Value = 17.0843614843986

$ws.Range("B20").Value = 8.047998458674064
$ws.Range("D20").Value = 3.997096436848766
$ws.Range("E20").Value = 19.812526468241
$ws.Range("F20").Value = 19.51083991699429
$ws.Range("G20").Value = 20.94309372373494
$ws.Range("H20").Value = 11.87925516577478
$ws.Range("K20").Value = 12.92163090230019
$ws.Range("N20").Value = 16.98856050114636
$ws.Range("O20").Value = 17.08251776923413

$ws.Range("B21").Value = 8.161001645408783
$ws.Range("D21").Value = 4.087032179491104
$ws.Range("E21").Value = 21.0161149778842
$ws.Range("F21").Value = 19.62644485102844
$ws.Range("G21").Value = 21.18404481958074
$ws.Range("H21").Value = 11.8499235063415
$ws.Range("K21").Value = 13.42995008626682
$ws.Range("N21").Value = 16.93645036261964
$ws.Range("O21").Value = 17.08202561087586

$ws.Range("B22").Value = 8.23544071614792
$ws.Range("D22").Value = 4.144268540075872
$ws.Range("E22").Value = 21.76730853554939
$ws.Range("F22").Value = 19.70716103304276
$ws.Range("G22").Value = 21.34790057866475
$ws.Range("H22").Value = 11.83280379749725
$ws.Range("K22").Value = 13.75043375987534
$ws.Range("N22").Value = 16.90391896822141
$ws.Range("O22").Value = 17.08595970486436

$ws.Range("B23").Value = 8.195665396153029
$ws.Range("D23").Value = 4.113864764761852
$ws.Range("E23").Value = 21.36961256187053
$ws.Range("F23").Value = 19.66361511818212
$ws.Range("G23").Value = 21.25988104800857
$ws.Range("H23").Value = 11.84175301850421
$ws.Range("K23").Value = 13.58047478658009
$ws.Range("N23").Value = 16.92114223148145
$ws.Range("O23").Value = 17.08347034429667

$ws.Range("B24").Value = 8.046163087725397
$ws.Range("D24").Value = 3.995601206290842
$ws.Range("E24").Value = 19.79225402197298
$ws.Range("F24").Value = 19.50903886954215
$ws.Range("G24").Value = 20.93926696796733
$ws.Range("H24").Value = 11.87976912130356
$ws.Range("K24").Value = 12.91312834871081
$ws.Range("N24").Value = 16.98943802319538
$ws.Range("O24").Value = 17.08259758274772

$ws.Range("B25").Value = 7.88819599215838
$ws.Range("D25").Value = 3.861767063588252
$ws.Range("E25").Value = 17.93763362085426
$ws.Range("F25").Value = 19.365029617229
$ws.Range("G25").Value = 20.62247266412152
$ws.Range("H25").Value = 11.92954672732692
$ws.Range("K25").Value = 12.14472845839613
$ws.Range("N25").Value = 17.06968075015475
$ws.Range("O25").Value = 17.09988314962726
